# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Tue Sep 24 04:13:52 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.071.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.623.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.624.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.095.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.949.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.620.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "339.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  -4.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "536.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0801"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "169.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.65%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "168.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.622"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0241"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0959"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
